# Update activity stats: reorder per-match rows (runs/balls/fours/sixes)
# for Eoin Morgan vs Kolkata Knight Riders down to the values now reflected
# in the Excel form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these as text values (matching the existing t="str" cell storage)
# by forcing a text number format before assigning the new value.
$updates = @{
    "C2" = "7";  "D2" = "10"; "E2" = "1"; "F2" = "0";
    "C3" = "24"; "D3" = "23"; "E3" = "2"; "F3" = "1";
    "C4" = "8";  "D4" = "12";             "F4" = "0";
    "C5" = "42"; "D5" = "29"; "E5" = "3"; "F5" = "2";
    "C7" = "44"; "D7" = "18";             "F7" = "5";
    "C8" = "34";             "E8" = "1"; "F8" = "2";
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
